$wb = $excel.ActiveWorkbook

# --- BGDPbES sheet: restore formula in B12 pulling from Wgtd Avg Expected Cap Factors ---
$wsBgd = $wb.Worksheets.Item("BGDPbES")
$wsBgd.Range("B12").Formula = "='Wgtd Avg Expected Cap Factors'!B12"

# Select BGDPbES sheet momentarily to set its selection to B10:AK11
$wsBgd.Select() | Out-Null
$wsBgd.Range("B10:AK11").Select() | Out-Null

# --- About sheet: update selection & re-activate as the selected tab ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Select() | Out-Null
$wsAbout.Range("A30").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
